# Update countries & provincias Spain
# Refresh the COVID-19 country data table in-place: update the "last
# updated" timestamp and rewrite every country row (name + stats) with
# the latest figures. Row positions/styles are preserved; only values
# change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Updated timestamp shown above the table
$ws.Range("A1").Value = 'Datos actualizados a 20 de Marzo de 2020 a las 05:46'

# Full country data: row, country, CasosTotales, NuevosCasos, CasosActivos,
# Recuperados, CasosCriticos, MuertesHoy, Muertes
$data = @(
    @(4, 'China', 80967, 39, 71150, 6569, 2136, 3, 3248),
    @(5, 'Italia', 41035, 0, 4440, 33190, 2498, 0, 3405),
    @(6, 'Iran', 18407, 0, 5979, 11144, 0, 0, 1284),
    @(7, 'España', 18077, 0, 1107, 16139, 939, 0, 831),
    @(8, 'Alemania', 15320, 0, 115, 15161, 2, 0, 44),
    @(9, 'Estados Unidos', 14340, 551, 125, 13998, 64, 9, 217),
    @(10, 'Francia', 10995, 0, 1295, 9328, 1122, 0, 372),
    @(11, 'Corea del Sur', 8652, 87, 2233, 6325, 59, 3, 94),
    @(12, 'Suiza', 4222, 0, 15, 4164, 0, 0, 43),
    @(13, 'Reino Unido', 3269, 0, 65, 3060, 20, 0, 144),
    @(14, 'Paises Bajos', 2460, 0, 2, 2382, 45, 0, 76),
    @(15, 'Austria', 2196, 17, 9, 2181, 13, 0, 6),
    @(16, 'Belgica', 1795, 0, 165, 1609, 130, 0, 21),
    @(17, 'Noruega', 1790, 0, 1, 1782, 27, 0, 7),
    @(18, 'Suecia', 1439, 0, 16, 1412, 21, 0, 11),
    @(19, 'Dinamarca', 1151, 0, 1, 1144, 30, 0, 6),
    @(20, 'Japon', 943, 0, 191, 719, 46, 0, 33),
    @(21, 'Malasia', 900, 0, 75, 823, 15, 0, 2),
    @(22, 'Canada', 873, 0, 11, 850, 1, 0, 12),
    @(23, 'Portugal', 786, 0, 4, 778, 20, 0, 4),
    @(24, 'Australia', 785, 29, 46, 732, 1, 0, 7),
    @(25, 'Crucero', 712, 0, 527, 178, 14, 0, 7),
    @(26, 'Chequia', 694, 0, 3, 691, 6, 0, 0),
    @(27, 'Israel', 677, 0, 14, 663, 6, 0, 0),
    @(28, 'Brasil', 647, 7, 2, 638, 18, 0, 7),
    @(29, 'Irlanda', 557, 0, 5, 549, 6, 0, 3),
    @(30, 'Grecia', 464, 0, 19, 439, 16, 0, 6),
    @(31, 'Catar', 460, 0, 10, 450, 6, 0, 0),
    @(32, 'Pakistan', 454, 0, 13, 439, 0, 0, 2),
    @(33, 'Finlandia', 400, 0, 10, 390, 2, 0, 0),
    @(34, 'Turquia', 359, 0, 0, 355, 0, 0, 4),
    @(35, 'Polonia', 355, 0, 13, 337, 3, 0, 5),
    @(36, 'Singapur', 345, 0, 124, 221, 14, 0, 0),
    @(37, 'Chile', 342, 0, 0, 342, 6, 0, 0),
    @(38, 'Luxemburgo', 335, 0, 6, 325, 1, 0, 4),
    @(39, 'Islandia', 330, 0, 5, 325, 1, 0, 0),
    @(40, 'Eslovenia', 319, 0, 0, 318, 6, 0, 1),
    @(41, 'Indonesia', 309, 0, 15, 269, 0, 0, 25),
    @(42, 'Barein', 278, 0, 110, 167, 4, 0, 1),
    @(43, 'Rumania', 277, 0, 25, 252, 5, 0, 0),
    @(44, 'Arabia Saudita', 274, 0, 8, 266, 0, 0, 0),
    @(45, 'Tailandia', 272, 0, 42, 229, 1, 0, 1),
    @(46, 'Estonia', 267, 0, 1, 266, 1, 0, 0),
    @(47, 'Ecuador', 260, 0, 1, 256, 2, 0, 3),
    @(48, 'Egipto', 256, 0, 42, 207, 0, 0, 7),
    @(49, 'Peru', 234, 0, 1, 230, 7, 2, 3),
    @(50, 'Filipinas', 217, 0, 8, 192, 1, 0, 17),
    @(51, 'Hong Kong', 208, 0, 98, 106, 4, 0, 4),
    @(52, 'Rusia', 199, 0, 8, 190, 0, 0, 1),
    @(53, 'India', 195, 1, 20, 171, 0, 0, 4),
    @(54, 'Irak', 192, 0, 49, 130, 0, 0, 13),
    @(55, 'Mexico', 164, 46, 4, 159, 1, 0, 1),
    @(56, 'Libano', 157, 0, 4, 149, 3, 0, 4),
    @(57, 'Sudafrica', 150, 0, 0, 150, 0, 0, 0),
    @(58, 'Kuwait', 148, 0, 18, 130, 5, 0, 0),
    @(59, 'San Marino', 144, 0, 4, 126, 12, 0, 14),
    @(60, 'Emiratos Arabes Unidos', 140, 0, 31, 109, 2, 0, 0),
    @(61, 'Panama', 137, 0, 1, 135, 7, 0, 1),
    @(62, 'Colombia', 128, 20, 1, 127, 0, 0, 0),
    @(63, 'Argentina', 128, 0, 3, 122, 0, 0, 3),
    @(64, 'Eslovaquia', 124, 0, 0, 124, 2, 0, 0),
    @(65, 'Armenia', 122, 0, 1, 121, 2, 0, 0),
    @(66, 'Croacia', 110, 0, 5, 104, 0, 0, 1),
    @(67, 'Taiwan', 108, 0, 26, 81, 0, 0, 1),
    @(68, 'Bulgaria', 107, 0, 0, 104, 0, 0, 3),
    @(69, 'Serbia', 103, 0, 1, 102, 4, 0, 0),
    @(70, 'Uruguay', 94, 15, 0, 94, 0, 0, 0),
    @(71, 'Argelia', 90, 0, 32, 49, 0, 0, 9),
    @(72, 'Costa Rica', 89, 2, 0, 88, 2, 0, 1),
    @(73, 'Letonia', 86, 0, 1, 85, 0, 0, 0),
    @(74, 'Vietnam', 85, 0, 16, 69, 0, 0, 0),
    @(75, 'Principado de Andorra', 74, 0, 1, 73, 0, 0, 0),
    @(76, 'Brunei', 73, 0, 0, 73, 2, 0, 0),
    @(77, 'Hungria', 73, 0, 2, 70, 4, 0, 1),
    @(78, 'Islas Feroe', 72, 0, 1, 71, 0, 0, 0),
    @(79, 'Jordania', 69, 0, 1, 68, 0, 0, 0),
    @(80, 'Republica de Chipre', 67, 0, 0, 67, 1, 0, 0),
    @(81, 'Bosnia y Herzegovina', 64, 0, 2, 62, 0, 0, 0),
    @(82, 'Albania', 64, 0, 0, 62, 2, 0, 2),
    @(83, 'Marruecos', 63, 0, 2, 59, 1, 0, 2),
    @(84, 'Sri Lanka', 60, 0, 3, 57, 0, 0, 0),
    @(85, 'Malta', 53, 0, 2, 51, 0, 0, 0),
    @(86, 'Bielorrusia', 51, 0, 5, 46, 0, 0, 0),
    @(87, 'Republica de Macedonia', 50, 0, 1, 49, 1, 0, 0),
    @(88, 'Kazajistan', 49, 5, 0, 49, 0, 0, 0),
    @(89, 'Moldavia', 49, 0, 1, 47, 3, 0, 1),
    @(90, 'Lituania', 48, 0, 1, 47, 1, 0, 0),
    @(91, 'Oman', 48, 0, 13, 35, 0, 0, 0),
    @(92, 'Estado de Palestina', 47, 0, 0, 47, 0, 0, 0),
    @(93, 'Azerbaiyan', 44, 0, 7, 36, 0, 0, 1),
    @(94, 'Venezuela', 42, 0, 0, 42, 0, 0, 0),
    @(95, 'Georgia', 40, 0, 1, 39, 1, 0, 0),
    @(96, 'Nueva Zelanda', 39, 11, 0, 39, 0, 0, 0),
    @(97, 'Tunez', 39, 0, 1, 37, 2, 0, 1),
    @(98, 'Camboya', 37, 0, 1, 36, 0, 0, 0),
    @(99, 'Senegal', 36, 0, 2, 34, 0, 0, 0),
    @(100, 'Republica Dominicana', 34, 0, 0, 32, 0, 0, 2),
    @(101, 'Guadalupe', 33, 0, 0, 33, 0, 0, 0),
    @(102, 'Burkina Faso', 33, 0, 0, 32, 0, 0, 1),
    @(103, 'Liechtenstein', 28, 0, 0, 28, 0, 0, 0),
    @(104, 'Reunion', 28, 0, 0, 28, 0, 0, 0),
    @(105, 'Ucrania', 26, 0, 0, 23, 0, 0, 3),
    @(106, 'Honduras', 24, 12, 0, 24, 0, 0, 0),
    @(107, 'Uzbekistan', 23, 0, 0, 23, 0, 0, 0),
    @(108, 'Martinica', 23, 0, 0, 22, 0, 0, 1),
    @(109, 'Afganistan', 22, 0, 1, 21, 0, 0, 0),
    @(110, 'Banglades', 18, 0, 3, 14, 0, 0, 1),
    @(111, 'Bolivia', 17, 2, 0, 17, 0, 0, 0),
    @(112, 'Macao', 17, 0, 10, 7, 0, 0, 0),
    @(113, 'Cuba', 16, 5, 0, 15, 0, 0, 1),
    @(114, 'Guayana Francesa', 15, 0, 0, 15, 0, 0, 0),
    @(115, 'Jamaica', 15, 0, 2, 12, 0, 0, 1),
    @(116, 'Consejo Danes para los Refugiados', 14, 0, 0, 14, 0, 0, 0),
    @(117, 'Maldivas', 13, 0, 0, 13, 0, 0, 0),
    @(118, 'Camerun', 13, 0, 0, 13, 0, 0, 0),
    @(119, 'Montenegro', 13, 0, 0, 13, 0, 0, 0),
    @(120, 'Paraguay', 13, 0, 0, 13, 1, 0, 0),
    @(121, 'Guam', 12, 0, 0, 12, 0, 0, 0),
    @(122, 'Nigeria', 12, 0, 1, 11, 0, 0, 0),
    @(123, 'Ghana', 11, 0, 0, 11, 0, 0, 0),
    @(124, 'Ruanda', 11, 0, 0, 11, 0, 0, 0),
    @(125, 'Monaco', 10, 0, 0, 10, 0, 0, 0),
    @(126, 'Gibraltar', 10, 0, 2, 8, 0, 0, 0),
    @(127, 'Trinidad yTobago', 9, 0, 0, 9, 0, 0, 0),
    @(128, 'Etiopia', 9, 2, 0, 9, 0, 0, 0),
    @(129, 'Guatemala', 9, 0, 0, 8, 0, 0, 1),
    @(130, 'Costa de Marfil', 9, 0, 1, 8, 0, 0, 0),
    @(131, 'Mauricio', 7, 0, 0, 7, 0, 0, 0),
    @(132, 'Kenia', 7, 0, 0, 7, 0, 0, 0),
    @(133, 'Guinea Ecuatorial', 6, 0, 0, 6, 0, 0, 0),
    @(134, 'Mongolia', 6, 0, 0, 6, 0, 0, 0),
    @(135, 'Polinesia Francesa', 6, 0, 0, 6, 0, 0, 0),
    @(136, 'Tanzania', 6, 0, 0, 6, 0, 0, 0),
    @(137, 'Kirguistan', 6, 3, 0, 6, 0, 0, 0),
    @(138, 'Seychelles', 6, 0, 0, 6, 0, 0, 0),
    @(139, 'Puerto Rico', 6, 0, 0, 6, 0, 0, 0),
    @(140, 'Barbados', 5, 0, 0, 5, 0, 0, 0),
    @(141, 'Guyana', 5, 0, 0, 4, 0, 0, 1),
    @(142, 'Aruba', 5, 0, 1, 4, 0, 0, 0),
    @(143, 'Mayotte', 4, 0, 0, 4, 0, 0, 0),
    @(144, 'Islas Virgenes de los Estados Unidos', 3, 0, 0, 3, 0, 0, 0),
    @(145, 'Congo', 3, 0, 0, 3, 0, 0, 0),
    @(146, 'Gabon', 3, 0, 0, 3, 0, 0, 0),
    @(147, 'San Bartolome', 3, 0, 0, 3, 0, 0, 0),
    @(148, 'Bahamas', 3, 0, 0, 3, 0, 0, 0),
    @(149, 'Namibia', 3, 0, 0, 3, 0, 0, 0),
    @(150, 'San Martin (Parte Francesa)', 3, 0, 0, 3, 0, 0, 0),
    @(151, 'Islas Caimanes', 3, 0, 0, 2, 0, 0, 1),
    @(152, 'Curazao', 3, 0, 0, 2, 0, 0, 1),
    @(153, 'Nueva Caledonia', 2, 0, 0, 2, 0, 0, 0),
    @(154, 'Mauritania', 2, 0, 0, 2, 0, 0, 0),
    @(155, 'Zambia', 2, 0, 0, 2, 0, 0, 0),
    @(156, 'Bermudas', 2, 0, 0, 2, 0, 0, 0),
    @(157, 'Groenlandia', 2, 0, 0, 2, 0, 0, 0),
    @(158, 'Haiti', 2, 2, 0, 2, 0, 0, 0),
    @(159, 'Benin', 2, 0, 0, 2, 0, 0, 0),
    @(160, 'Liberia', 2, 0, 0, 2, 0, 0, 0),
    @(161, 'Santa Lucia', 2, 0, 0, 2, 0, 0, 0),
    @(162, 'Sudan', 2, 0, 0, 1, 0, 0, 1),
    @(163, 'Montserrat', 1, 0, 0, 1, 0, 0, 0),
    @(164, 'Gambia', 1, 0, 0, 1, 0, 0, 0),
    @(165, 'Niger', 1, 0, 0, 1, 0, 0, 0),
    @(166, 'Guinea', 1, 0, 0, 1, 0, 0, 0),
    @(167, 'Fiyi', 1, 0, 0, 1, 0, 0, 0),
    @(168, 'Antigua y Barbuda', 1, 0, 0, 1, 0, 0, 0),
    @(169, 'Surinam', 1, 0, 0, 1, 0, 0, 0),
    @(170, 'Nicaragua', 1, 0, 0, 1, 0, 0, 0),
    @(171, 'San Martin (Parte Holandesa)', 1, 0, 0, 1, 0, 0, 0),
    @(172, 'Republica de Africa Central', 1, 0, 0, 1, 0, 0, 0),
    @(173, 'Santa Sede', 1, 0, 0, 1, 0, 0, 0),
    @(174, 'Butan', 1, 0, 0, 1, 0, 0, 0),
    @(175, 'San Vicente y las Granadinas', 1, 0, 0, 1, 0, 0, 0),
    @(176, 'Somalia', 1, 0, 0, 1, 0, 0, 0),
    @(177, 'Suazilandia', 1, 0, 0, 1, 0, 0, 0),
    @(178, 'Isla de Man', 1, 0, 0, 1, 0, 0, 0),
    @(179, 'Republica de Yibuti', 1, 0, 0, 1, 0, 0, 0),
    @(180, 'Republica del Chad', 1, 0, 0, 1, 0, 0, 0),
    @(181, 'Togo', 1, 0, 0, 1, 0, 0, 0),
    @(182, 'El Salvador', 1, 0, 0, 1, 0, 0, 0),
    @(183, 'Nepal', 1, 0, 1, 0, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
